$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 606.07526
$ws.Range("J17").Value = 606.07526
$ws.Range("L17").Value = 1818.22578
$ws.Range("N17").Value = -2154.22578
$ws.Range("H32").Value = 525.2222
$ws.Range("I32").Value = 522.475
$ws.Range("K32").Value = 522.475
$ws.Range("M32").Value = -196.475
$ws.Range("H40").Value = 1198.9744
$ws.Range("I40").Value = 1252.5416
$ws.Range("J40").Value = 1113.2667
$ws.Range("K40").Value = 1252.5416
$ws.Range("L40").Value = 1113.2667
$ws.Range("M40").Value = -1077.5416
$ws.Range("N40").Value = -1463.2667
$ws.Range("H53").Value = 191
$ws.Range("I53").Value = 156.3125
$ws.Range("J53").Value = 252.66667
$ws.Range("K53").Value = 156.3125
$ws.Range("L53").Value = 252.66667
$ws.Range("M53").Value = 480.6875
$ws.Range("N53").Value = -1526.66667
$ws.Range("H61").Value = 907.5
$ws.Range("I61").Value = 907.5
$ws.Range("K61").Value = 2722.5
$ws.Range("M61").Value = -2550.5
$ws.Range("H64").Value = 3017
$ws.Range("I64").Value = 2999.8
$ws.Range("K64").Value = 2999.8
$ws.Range("M64").Value = -2751.8
$ws.Range("H67").Value = 3017
$ws.Range("I67").Value = 2999.8
$ws.Range("K67").Value = 2999.8
$ws.Range("M67").Value = -2141.8
$ws.Range("H70").Value = 1000
$ws.Range("J70").Value = 1000
$ws.Range("L70").Value = 3000
$ws.Range("N70").Value = -3540
$ws.Range("H73").Value = 1000
$ws.Range("J73").Value = 1000
$ws.Range("L73").Value = 3000
$ws.Range("N73").Value = -4872
$ws.Range("H132").Value = 1856593.5
$ws.Range("I132").Value = 2067052
$ws.Range("J132").Value = 4559.2
$ws.Range("K132").Value = 6201156
$ws.Range("L132").Value = 13677.6
$ws.Range("M132").Value = -6198626
$ws.Range("N132").Value = -18737.6
$ws.Range("H135").Value = 892.449
$ws.Range("I135").Value = 193.05556
$ws.Range("K135").Value = 1737.50004
$ws.Range("M135").Value = 797.4999599999999
$ws.Range("H137").Value = 1270.2162
$ws.Range("I137").Value = 1003.1667
$ws.Range("J137").Value = 2414.7144
$ws.Range("K137").Value = 3009.5001
$ws.Range("L137").Value = 7244.1432
$ws.Range("M137").Value = -459.5001000000002
$ws.Range("N137").Value = -12344.1432
$ws.Range("H138").Value = 12052021
$ws.Range("I138").Value = 28572324
$ws.Range("J138").Value = 5967.9165
$ws.Range("K138").Value = 85716972
$ws.Range("L138").Value = 17903.7495
$ws.Range("M138").Value = -85711832
$ws.Range("N138").Value = -28183.7495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7739.3076
$ws.Range("I32").Value = 4625.2095
$ws.Range("J32").Value = 19806.438
$ws.Range("K32").Value = 4625.2095
$ws.Range("L32").Value = 19806.438
$ws.Range("M32").Value = -4338.2095
$ws.Range("N32").Value = -20380.438
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H63").Value = 2427.7964
$ws.Range("I63").Value = 2449.8958
$ws.Range("J63").Value = 2251
$ws.Range("K63").Value = 2449.8958
$ws.Range("L63").Value = 2251
$ws.Range("M63").Value = -1763.8958
$ws.Range("N63").Value = -3623
$ws.Range("H66").Value = 2427.7964
$ws.Range("I66").Value = 2449.8958
$ws.Range("J66").Value = 2251
$ws.Range("K66").Value = 12249.479
$ws.Range("L66").Value = 11255
$ws.Range("M66").Value = -8817.478999999999
$ws.Range("N66").Value = -18119
$ws.Range("H74").Value = 895.62067
$ws.Range("I74").Value = 822.4400000000001
$ws.Range("J74").Value = 1353
$ws.Range("K74").Value = 822.4400000000001
$ws.Range("L74").Value = 1353
$ws.Range("M74").Value = 51.55999999999995
$ws.Range("N74").Value = -3101
$ws.Range("H77").Value = 895.62067
$ws.Range("I77").Value = 822.4400000000001
$ws.Range("J77").Value = 1353
$ws.Range("K77").Value = 4112.200000000001
$ws.Range("L77").Value = 6765
$ws.Range("M77").Value = 255.7999999999993
$ws.Range("N77").Value = -15501
$ws.Range("H122").Value = 2351.0256
$ws.Range("I122").Value = 1920.375
$ws.Range("J122").Value = 2650.6086
$ws.Range("K122").Value = 5761.125
$ws.Range("L122").Value = 7951.825800000001
$ws.Range("M122").Value = -3311.125
$ws.Range("N122").Value = -12851.8258
$ws.Range("H132").Value = 1995.0435
$ws.Range("I132").Value = 1153.7333
$ws.Range("J132").Value = 3572.5
$ws.Range("K132").Value = 3461.199900000001
$ws.Range("L132").Value = 10717.5
$ws.Range("M132").Value = -931.1999000000005
$ws.Range("N132").Value = -15777.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 12236
$ws.Range("I54").Value = 1845
$ws.Range("J54").Value = 53800
$ws.Range("K54").Value = 1845
$ws.Range("L54").Value = 53800
$ws.Range("M54").Value = -1361
$ws.Range("N54").Value = -54768
$ws.Range("H86").Value = 1550
$ws.Range("I86").Value = 1250
$ws.Range("J86").Value = 1850
$ws.Range("K86").Value = 1250
$ws.Range("L86").Value = 1850
$ws.Range("M86").Value = -127
$ws.Range("N86").Value = -4096
$ws.Range("H89").Value = 1550
$ws.Range("I89").Value = 1250
$ws.Range("J89").Value = 1850
$ws.Range("K89").Value = 6250
$ws.Range("L89").Value = 9250
$ws.Range("M89").Value = -634
$ws.Range("N89").Value = -20482
$ws.Range("H134").Value = 1218.64
$ws.Range("I134").Value = 894.94116
$ws.Range("J134").Value = 1906.5
$ws.Range("K134").Value = 2684.82348
$ws.Range("L134").Value = 5719.5
$ws.Range("M134").Value = -149.82348
$ws.Range("N134").Value = -10789.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1896.8036
$ws.Range("I31").Value = 1242.7021
$ws.Range("J31").Value = 5312.6665
$ws.Range("K31").Value = 1242.7021
$ws.Range("L31").Value = 5312.6665
$ws.Range("M31").Value = -947.7021
$ws.Range("N31").Value = -5902.6665
$ws.Range("H34").Value = 1896.8036
$ws.Range("I34").Value = 1242.7021
$ws.Range("J34").Value = 5312.6665
$ws.Range("K34").Value = 1242.7021
$ws.Range("L34").Value = 5312.6665
$ws.Range("M34").Value = -1040.7021
$ws.Range("N34").Value = -5716.6665
$ws.Range("H58").Value = 1116.2778
$ws.Range("I58").Value = 944.11536
$ws.Range("J58").Value = 1563.9
$ws.Range("K58").Value = 944.11536
$ws.Range("L58").Value = 1563.9
$ws.Range("M58").Value = -741.11536
$ws.Range("N58").Value = -1969.9
$ws.Range("H132").Value = 1556.4762
$ws.Range("I132").Value = 1397.2667
$ws.Range("J132").Value = 1954.5
$ws.Range("K132").Value = 4191.800099999999
$ws.Range("L132").Value = 5863.5
$ws.Range("M132").Value = -1661.800099999999
$ws.Range("N132").Value = -10923.5
$ws.Range("H134").Value = 1121.0834
$ws.Range("I134").Value = 958.34784
$ws.Range("J134").Value = 1655.7858
$ws.Range("K134").Value = 2875.04352
$ws.Range("L134").Value = 4967.357400000001
$ws.Range("M134").Value = -340.0435200000002
$ws.Range("N134").Value = -10037.3574
$ws.Range("H136").Value = 1116.2778
$ws.Range("I136").Value = 944.11536
$ws.Range("J136").Value = 1563.9
$ws.Range("K136").Value = 2832.34608
$ws.Range("L136").Value = 4691.700000000001
$ws.Range("M136").Value = -282.3460800000003
$ws.Range("N136").Value = -9791.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 875.871
$ws.Range("I131").Value = 453
$ws.Range("J131").Value = 999.2083
$ws.Range("K131").Value = 1359
$ws.Range("L131").Value = 2997.6249
$ws.Range("M131").Value = 3681
$ws.Range("N131").Value = -13077.6249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2212
$ws.Range("I80").Value = 2244.2222
$ws.Range("J80").Value = 2175.75
$ws.Range("K80").Value = 2244.2222
$ws.Range("L80").Value = 2175.75
$ws.Range("M80").Value = -1246.2222
$ws.Range("N80").Value = -4171.75
$ws.Range("H83").Value = 2212
$ws.Range("I83").Value = 2244.2222
$ws.Range("J83").Value = 2175.75
$ws.Range("K83").Value = 11221.111
$ws.Range("L83").Value = 10878.75
$ws.Range("M83").Value = -6229.111000000001
$ws.Range("N83").Value = -20862.75
$ws.Range("H102").Value = 2802.2693
$ws.Range("I102").Value = 1750.5333
$ws.Range("J102").Value = 4236.4546
$ws.Range("K102").Value = 1750.5333
$ws.Range("L102").Value = 4236.4546
$ws.Range("M102").Value = -128.5333000000001
$ws.Range("N102").Value = -7480.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2885.3333
$ws.Range("I7").Value = 2183.6
$ws.Range("J7").Value = 3762.5
$ws.Range("K7").Value = 2183.6
$ws.Range("L7").Value = 3762.5
$ws.Range("M7").Value = -2071.6
$ws.Range("N7").Value = -3986.5
$ws.Range("H82").Value = 2109
$ws.Range("I82").Value = 1500
$ws.Range("J82").Value = 2616.5
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 2616.5
$ws.Range("M82").Value = -1139
$ws.Range("N82").Value = -3338.5
$ws.Range("H85").Value = 2109
$ws.Range("I85").Value = 1500
$ws.Range("J85").Value = 2616.5
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 2616.5
$ws.Range("M85").Value = -252
$ws.Range("N85").Value = -5112.5
$ws.Range("H126").Value = 2885.3333
$ws.Range("I126").Value = 2183.6
$ws.Range("J126").Value = 3762.5
$ws.Range("K126").Value = 6550.799999999999
$ws.Range("L126").Value = 11287.5
$ws.Range("M126").Value = -4080.799999999999
$ws.Range("N126").Value = -16227.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 30000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -36240
$ws.Range("H128").Value = 37128.75
$ws.Range("J128").Value = 37128.75
$ws.Range("L128").Value = 37128.75
$ws.Range("N128").Value = -47088.75
$ws.Range("H132").Value = 1084.2766
$ws.Range("I132").Value = 816.2143
$ws.Range("J132").Value = 1479.3158
$ws.Range("K132").Value = 2448.6429
$ws.Range("L132").Value = 4437.9474
$ws.Range("M132").Value = 81.35710000000017
$ws.Range("N132").Value = -9497.947400000001
